# Update gh-pages to output generated at 456a3b4
# Refreshes the "想去人数" (interest count) figures in column F for the
# "展览" and "全部类型" sheets. All other cells are left untouched.

$wb = $excel.ActiveWorkbook

# -- Sheet "展览" (sheet1) --------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7732
$ws1.Range("F3").Value  = 7538
$ws1.Range("F4").Value  = 107
$ws1.Range("F5").Value  = 186
$ws1.Range("F8").Value  = 122
$ws1.Range("F10").Value = 143
$ws1.Range("F11").Value = 216
$ws1.Range("F12").Value = 103
$ws1.Range("F14").Value = 968
$ws1.Range("F16").Value = 39
$ws1.Range("F17").Value = 5
$ws1.Range("F19").Value = 88

# -- Sheet "全部类型" (sheet4) ------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 7732
$ws4.Range("F3").Value  = 7538
$ws4.Range("F4").Value  = 107
$ws4.Range("F5").Value  = 186
$ws4.Range("F8").Value  = 122
$ws4.Range("F10").Value = 143
$ws4.Range("F11").Value = 216
$ws4.Range("F12").Value = 103
$ws4.Range("F14").Value = 969
$ws4.Range("F16").Value = 39
$ws4.Range("F17").Value = 5
$ws4.Range("F19").Value = 88
